# Update "Out of PO.xlsx" player roster data (A2:C19) to new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Grant Williams", "PF,C", "Charlotte Hornets"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Ja Morant", "PG", "Memphis Grizzlies")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
